$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.407.03"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.518.26"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.25"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.73"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "2.523.65"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.359"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "2.964.89"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.54"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "59.301.09"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "2.518.00"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.14"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.69"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.20"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.427"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.99"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").Value = "0.0₃0780"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.98"
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.46"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  -7.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.53"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.96"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  -6.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.50"
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.599"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.85"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.16"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0225"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.88"
$ws.Range("E51").Value = "  -1.84%  "
